# Daily attendance processing - 2025-12-08 19:25:29
# Re-sort the "Recorded By" (column G) comma-separated list of names/emails
# in case-insensitive alphabetical order for every data row on the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $val = $cell.Value()

    if ($val -ne $null -and $val -ne "") {
        $parts = $val -split ","
        $trimmed = @()
        foreach ($p in $parts) {
            $trimmed += $p.Trim()
        }

        $sorted = $trimmed | Sort-Object { $_.ToLower() }

        $newVal = [string]::Join(", ", $sorted)

        if ($newVal -ne $val) {
            $cell.Value = $newVal
        }
    }
}
